$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already had rows 3-11 of slang entries, but rows 4-11 were
# missing their "id_giria" (col A) and "dificuldade" (col G) values.
# Fill those in now that the data set has been fleshed out.
$idValues = @{ 4 = 3; 5 = 4; 6 = 5; 7 = 6; 8 = 7; 9 = 8; 10 = 9; 11 = 10 }
$gValues  = @{ 4 = 1; 5 = 1; 6 = 2; 7 = 1; 8 = 1; 9 = 2; 10 = 2; 11 = 1 }

foreach ($row in 4..11) {
    $ws.Cells.Item($row, 1).Value = $idValues[$row]
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# Move the view/selection: was scrolled to column B with C13 selected,
# now scrolled to column C with G12 selected.
$window = $excel.ActiveWindow
$window.ScrollColumn = 3
$window.ScrollRow = 1
$ws.Range("G12").Select()
